$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# "Manage Error Types" chapter - extend the final sentence of the chapter's
# intro paragraph so that:
#   "Using the Error Type Details form."
# becomes:
#   "Using the Error Type Details form, the user can add new error types."
# ---------------------------------------------------------------------------

$rng = $d.Content
[void]$rng.Find.Execute("Error Type Details form")
$rng.Collapse(0)
$insertPos = $rng.Start

$newText = ", the user can add new error types"
[void]$rng.InsertAfter($newText)

# Word always leaves the hidden "_GoBack" bookmark at the location of the
# most recent edit - here, right after "can", which is where the user's
# cursor was left before finishing the sentence with " add new error types".
$bmPos = $insertPos + 14
$bmRng = $d.Range($bmPos, $bmPos)
[void]$d.Bookmarks.Add("_GoBack", $bmRng)
